$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.533.05"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").Value = "2.029.45"
$ws.Range("E3").Value = "  +1.86%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'233.63"
$ws.Range("E5").Value = "  -8.59%  "
$ws.Range("D6").Value = "'0.602"
$ws.Range("E6").Value = "  -1.54%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'55.39"
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("D10").Value = "'57.32"
$ws.Range("E10").Value = "  +3.38%  "
$ws.Range("D11").Value = "'0.0750"
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("D12").Value = "'0.101"
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.326.67"
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'14.40"
$ws.Range("E14").Value = "  +1.81%  "
$ws.Range("D15").Value = "'20.22"
$ws.Range("E15").Value = "  -4.70%  "
$ws.Range("E16").Value = "  -3.38%  "
$ws.Range("D17").Value = "'5.11"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").Value = "2.012.30"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("D19").Value = "36.715.34"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").Value = "'67.85"
$ws.Range("E20").Value = "  -3.92%  "
$ws.Range("D21").Value = "'5.59"
$ws.Range("E21").Value = "  +10.84%  "
$ws.Range("E22").Value = "  -2.65%  "
$ws.Range("E23").Value = "  -5.81%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("E26").Value = "  -6.27%  "
$ws.Range("D27").Value = "'162.87"
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("D28").Value = "'8.68"
$ws.Range("E28").Value = "  -2.17%  "
$ws.Range("D29").Value = "'0.131"
$ws.Range("E29").Value = "  +5.84%  "
$ws.Range("D30").Value = "'19.00"
$ws.Range("E30").Value = "  -1.84%  "
$ws.Range("D31").Value = "'1.35"
$ws.Range("E31").Value = "  +1.73%  "
$ws.Range("E32").Value = "  -0.99%  "
$ws.Range("E33").Value = "  -2.97%  "
$ws.Range("D34").Value = "'0.0604"
$ws.Range("E34").Value = "  -3.95%  "
$ws.Range("E36").Value = "  -1.75%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("D39").Value = "'5.85"
$ws.Range("E39").Value = "  +6.23%  "
$ws.Range("D40").Value = "'3.29"
$ws.Range("E40").Value = "  -5.76%  "
$ws.Range("D41").Value = "'2.92"
$ws.Range("E41").Value = "  -2.75%  "
$ws.Range("D42").Value = "'0.0943"
$ws.Range("E42").Value = "  +3.38%  "
$ws.Range("D43").Value = "1.467.18"
$ws.Range("E43").Value = "  +1.47%  "
$ws.Range("D44").Value = "'94.63"
$ws.Range("E44").Value = "  +6.91%  "
$ws.Range("E45").Value = "  -1.82%  "
$ws.Range("E46").Value = "  -4.13%  "
$ws.Range("D47").Value = "'15.68"
$ws.Range("E47").Value = "  +1.28%  "
$ws.Range("E48").Value = "  +36.72%  "
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("E51").Value = "  +0.56%  "
